# Applies the "13-12-2022" crypto price-list refresh (GitHub Actions bot update):
# - column G ("Hora") bumps from 4 -> 5 for every data row
# - column D ("Price") refreshes to the latest quote on many rows
# - rows 14-26 and 42-43 shuffle their Coin/Link/Price/Volume columns
#   (coins drop out / move position in the ranking)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Force text storage so numeric-looking strings (e.g. "6.250", "0.1120")
    # keep their exact formatting instead of being coerced to a Double.
    $cell = $ws.Range($range)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

Set-TextCell "D2" '268.95'
Set-TextCell "G2" '5'

Set-TextCell "D3" '21.38'
Set-TextCell "G3" '5'

Set-TextCell "D4" '6.259'
Set-TextCell "G4" '5'

Set-TextCell "D5" '0.06209'
Set-TextCell "G5" '5'

Set-TextCell "D6" '3.571'
Set-TextCell "G6" '5'

Set-TextCell "D7" '6.538'
Set-TextCell "G7" '5'

Set-TextCell "D8" '1.396'
Set-TextCell "G8" '5'

Set-TextCell "D9" '0.8256'
Set-TextCell "G9" '5'

Set-TextCell "D10" '0.1638'
Set-TextCell "G10" '5'

Set-TextCell "D11" '0.08254'
Set-TextCell "G11" '5'

Set-TextCell "D12" '0.03558'
Set-TextCell "G12" '5'

Set-TextCell "D13" '0.03187'
Set-TextCell "G13" '5'

# row 14 -> BitMartToken
Set-TextCell "B14" 'BitMartToken'
Set-TextCell "C14" 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextCell "D14" '0.09199'
Set-TextCell "E14" '13BitMartTokenBMX'
Set-TextCell "G14" '5'

# row 15 -> MCDex
Set-TextCell "B15" 'MCDex'
Set-TextCell "C15" 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextCell "D15" '3.763'
Set-TextCell "E15" '14MCDexMCB'
Set-TextCell "G15" '5'

# row 16 -> BitForexToken
Set-TextCell "B16" 'BitForexToken'
Set-TextCell "C16" 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextCell "D16" '0.001638'
Set-TextCell "E16" '15BitForexTokenBF'
Set-TextCell "G16" '5'

# row 17 -> CoinExToken
Set-TextCell "B17" 'CoinExToken'
Set-TextCell "C17" 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextCell "D17" '0.04679'
Set-TextCell "E17" '16CoinExTokenCET'
Set-TextCell "G17" '5'

# row 18 -> TigerCash
Set-TextCell "B18" 'TigerCash'
Set-TextCell "C18" 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextCell "D18" '0.006463'
Set-TextCell "E18" '17TigerCashTCH'
Set-TextCell "G18" '5'

# row 19 -> HotbitToken
Set-TextCell "B19" 'HotbitToken'
Set-TextCell "C19" 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextCell "D19" '0.006191'
Set-TextCell "E19" '18HotbitTokenHTB'
Set-TextCell "G19" '5'

# row 20 -> BitKan
Set-TextCell "B20" 'BitKan'
Set-TextCell "C20" 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextCell "D20" '0.001069'
Set-TextCell "E20" '19BitKanKAN'
Set-TextCell "G20" '5'

# row 21 -> NitroEx
Set-TextCell "B21" 'NitroEx'
Set-TextCell "C21" 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextCell "D21" '0.0001501'
Set-TextCell "E21" '20NitroExNTX'
Set-TextCell "G21" '5'

# row 22 -> LEO
Set-TextCell "B22" 'LEO'
Set-TextCell "C22" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextCell "D22" '3.723'
Set-TextCell "E22" '21LEOLEO'
Set-TextCell "G22" '5'

# row 23 -> BTSEToken
Set-TextCell "B23" 'BTSEToken'
Set-TextCell "C23" 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextCell "D23" '2.286'
Set-TextCell "E23" '22BTSETokenBTSE'
Set-TextCell "G23" '5'

# row 24 -> One
Set-TextCell "B24" 'One'
Set-TextCell "C24" 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextCell "D24" '0.01364'
Set-TextCell "E24" '23OneONE'
Set-TextCell "G24" '5'

# row 25 -> BitpandaEcosystemToken
Set-TextCell "B25" 'BitpandaEcosystemToken'
Set-TextCell "C25" 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextCell "D25" '0.3318'
Set-TextCell "E25" '24BitpandaEcosystemTokenBEST'
Set-TextCell "G25" '5'

# row 26 -> ProBitToken
Set-TextCell "B26" 'ProBitToken'
Set-TextCell "C26" 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextCell "D26" '0.1242'
Set-TextCell "E26" '25ProBitTokenPROB'
Set-TextCell "G26" '5'

Set-TextCell "G27" '5'

Set-TextCell "D28" '0.0002714'
Set-TextCell "G28" '5'

Set-TextCell "G29" '5'

Set-TextCell "G30" '5'

Set-TextCell "G31" '5'

Set-TextCell "G32" '5'

Set-TextCell "G33" '5'

Set-TextCell "G34" '5'

Set-TextCell "G35" '5'

Set-TextCell "G36" '5'

Set-TextCell "G37" '5'

Set-TextCell "G38" '5'

Set-TextCell "G39" '5'

Set-TextCell "D40" '0.04714'
Set-TextCell "G40" '5'

Set-TextCell "D41" '0.006973'
Set-TextCell "G41" '5'

# row 42 -> CEJI
Set-TextCell "B42" 'CEJI'
Set-TextCell "C42" 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextCell "D42" '0.004002'
Set-TextCell "E42" '41CEJICEJI'
Set-TextCell "G42" '5'

# row 43 -> BKEXToken
Set-TextCell "B43" 'BKEXToken'
Set-TextCell "C43" 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextCell "D43" '0.1120'
Set-TextCell "E43" '42BKEXTokenBKK'
Set-TextCell "G43" '5'

Set-TextCell "D44" '0.01158'
Set-TextCell "G44" '5'

Set-TextCell "D45" '0.00006195'
Set-TextCell "G45" '5'

Set-TextCell "D46" '0.0009904'
Set-TextCell "G46" '5'

Set-TextCell "G47" '5'

Set-TextCell "D48" '0.8028'
Set-TextCell "G48" '5'

Set-TextCell "D49" '0.002339'
Set-TextCell "G49" '5'

Set-TextCell "D50" '0.00001901'
Set-TextCell "G50" '5'

Set-TextCell "G51" '5'
